$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.034.50'
$ws.Range('E2').Value = '  +0.50%  '
# Row 3
$ws.Range('D3').Value = '1.643.89'
$ws.Range('E3').Value = '  +0.63%  '
# Row 4
$ws.Range('E4').Value = '  +0.25%  '
# Row 5
$ws.Range('D5').Value = '216.39'
$ws.Range('E5').Value = '  +0.87%  '
# Row 6
$ws.Range('E6').Value = '  +0.68%  '
# Row 7
$ws.Range('E7').Value = '  +0.23%  '
# Row 8
$ws.Range('E8').Value = '  +0.60%  '
# Row 9
$ws.Range('E9').Value = '  +1.34%  '
# Row 10
$ws.Range('D10').Value = '19.65'
$ws.Range('E10').Value = '  +0.62%  '
# Row 11
$ws.Range('E11').Value = '  +0.66%  '
# Row 12
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.872.21'
$ws.Range('E12').Value = '  +0.64%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.28'
$ws.Range('E13').Value = '  +1.03%  '
# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.666.67'
$ws.Range('E14').Value = '  +1.71%  '
# Row 15
$ws.Range('E15').Value = '  +0.32%  '
# Row 16
$ws.Range('E16').Value = '  +1.44%  '
# Row 17
$ws.Range('D17').Value = '63.38'
$ws.Range('E17').Value = '  +0.96%  '
# Row 18
$ws.Range('D18').Value = '26.050.64'
$ws.Range('E18').Value = '  +0.56%  '
# Row 19
$ws.Range('E19').Value = '  +0.34%  '
# Row 20
$ws.Range('D20').Value = '195.56'
$ws.Range('E20').Value = '  +1.49%  '
# Row 21
$ws.Range('E21').Value = '  -0.41%  '
# Row 22
$ws.Range('D22').Value = '9.93'
$ws.Range('E22').Value = '  +0.01%  '
# Row 23
$ws.Range('D23').Value = '6.23'
$ws.Range('E23').Value = '  +0.00%  '
# Row 24
$ws.Range('E24').Value = '  +5.10%  '
# Row 25
$ws.Range('E25').Value = '  -0.14%  '
# Row 26
$ws.Range('E26').Value = '  +0.59%  '
# Row 27
$ws.Range('D27').Value = '143.51'
$ws.Range('E27').Value = '  +0.32%  '
# Row 28
$ws.Range('E28').Value = '  +0.90%  '
# Row 29
$ws.Range('D29').Value = '15.55'
$ws.Range('E29').Value = '  +0.74%  '
# Row 30
$ws.Range('E30').Value = '  +1.24%  '
# Row 31
$ws.Range('E31').Value = '  +0.18%  '
# Row 32
$ws.Range('E32').Value = '  -0.07%  '
# Row 33
$ws.Range('E33').Value = '  +1.51%  '
# Row 34
$ws.Range('D34').Value = '1.53'
$ws.Range('E34').Value = '  -2.70%  '
# Row 35
$ws.Range('E35').Value = '  +1.23%  '
# Row 36
$ws.Range('E36').Value = '  +0.73%  '
# Row 37
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.133.47'
$ws.Range('E37').Value = '  -0.35%  '
# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.541'
$ws.Range('E38').Value = '  -1.10%  '
# Row 39
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').Value = '  -0.61%  '
# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.0157'
$ws.Range('E40').Value = '  +0.54%  '
# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '5.50'
$ws.Range('E41').Value = '  +0.99%  '
# Row 42
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '99.33'
$ws.Range('E42').Value = '  +0.13%  '
# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.799'
$ws.Range('E43').Value = '  -0.78%  '
# Row 44
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.781.47'
$ws.Range('E44').Value = '  +0.67%  '
# Row 45
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0117'
$ws.Range('E45').Value = '  +4.03%  '
# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '56.82'
$ws.Range('E46').Value = '  +1.09%  '
# Row 47
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.0525'
$ws.Range('E47').Value = '  +0.51%  '
# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.47'
$ws.Range('E48').Value = '  +1.05%  '
# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.76'
$ws.Range('E49').Value = '  +1.89%  '
# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.414'
$ws.Range('E50').Value = '  -0.25%  '
# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0960'
$ws.Range('E51').Value = '  -0.01%  '
